$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right before the current row 539 (the oldest rows shift down),
# adding one more weekly batch of "Comercializadora del Agro de Limarí - Frutilla" quotes
# (commit: "Fruta / hortaliza, semanal").
$ws.Rows("539:541").Insert()

# Common values shared by the three new rows.
$mercadoId = 2
$mercado = "Comercializadora del Agro de Limarí"
$region = "Coquimbo"
$fecha = 45021
$codreg = 4
$tipo = "Fruta"
$productoId = 100101
$producto = "Berries"
$categoriaId = 100112025
$categoria = "Frutilla"
$variedad = "Sin especificar"
$unidad = "$/bandeja 7 kilos"
$origen = "Provincia de Melipilla"
$kgUnidad = 7

# Row 539 - Especial
$r = 539
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 500
$ws.Cells.Item($r, 14).Value = 13000
$ws.Cells.Item($r, 15).Value = 14000
$ws.Cells.Item($r, 16).Value = 13500
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 1929
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 540 - Primera
$r = 540
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 400
$ws.Cells.Item($r, 14).Value = 11000
$ws.Cells.Item($r, 15).Value = 12000
$ws.Cells.Item($r, 16).Value = 11500
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 1643
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 541 - Segunda
$r = 541
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 320
$ws.Cells.Item($r, 14).Value = 9000
$ws.Cells.Item($r, 15).Value = 10000
$ws.Cells.Item($r, 16).Value = 9500
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 1357
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Keep the date column's display style (yyyy-mm-dd hh:mm:ss) consistent with the rest of column D.
$ws.Range("D539:D541").NumberFormat = $ws.Range("D542").NumberFormat
